$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DEC-2020")

# --- 1. Move the legend block (rows 27-31) down to rows 36-40 ---
$ws.Range("C27").Copy($ws.Range("C36"))
$ws.Range("B28:C31").Copy($ws.Range("B37"))
$ws.Range("B27:C31").Clear()

# --- 2. Build new rows 20-21: empty 2-row merged placeholder (same look as rows 6/7, 13/14) ---
$ws.Range("A13:G14").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Range("A20").Value = 5
$ws.Range("B20").Value = 44184
$ws.Range("A21").Value = 6
$ws.Range("B21").Value = 44185
$ws.Range("C20:G21").Merge()

# --- 3. Build new rows 22-25 (normal data rows, same look as rows 15-18) ---
$ws.Range("A15:G18").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("A22").Value = 7
$ws.Range("B22").Value = 44186
$ws.Range("C22").Value = "QMVAR 2.0"
$ws.Range("D22").Value = "Design issue fixing"
$ws.Range("F22").Value = "WIP"

$ws.Range("A23").Value = 8
$ws.Range("B23").Value = 44187
$ws.Range("C23").Value = "QMVAR 2.0"
$ws.Range("D23").Value = "Design issue fixing"
$ws.Range("F23").Value = "WIP"

$ws.Range("A24").Value = 8
$ws.Range("B24").Value = 44188
$ws.Range("C24").Value = "QMVAR 2.0"
$ws.Range("D24").Value = "Design issue fixing"
$ws.Range("F24").Value = "Completed"

$ws.Range("A25").Value = 8
$ws.Range("B25").Value = 44189
$ws.Range("C25").Value = "QMVAR 2.0"
$ws.Range("D25").Value = "Design issue fixing in scroll bar"
$ws.Range("F25").Value = "Completed"

# --- 4. Row 26: single placeholder row (top edge of the new 3-row merged box) ---
$ws.Range("A15:G15").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Range("A26").Value = 8
$ws.Range("B26").Value = 44190
$ws.Range("C26:G26").ClearContents()
$ws.Range("C26:G26").Font.Size = 11
$ws.Range("C26:G26").WrapText = $false
$ws.Range("C26:G26").Borders.Item(7).LineStyle = 0
$ws.Range("C26:G26").Borders.Item(8).LineStyle = 0
$ws.Range("C26:G26").Borders.Item(9).LineStyle = 0
$ws.Range("C26:G26").Borders.Item(10).LineStyle = 0
$ws.Range("C26").Borders.Item(7).LineStyle = 1
$ws.Range("C26:G26").Borders.Item(8).LineStyle = 1
$ws.Range("G26").Borders.Item(10).LineStyle = 1

# --- 5. Rows 27-28: continuation of the 3-row merged box ---
$ws.Range("A27").Value = 5
$ws.Range("B27").Value = 44191
$ws.Range("A28").Value = 6
$ws.Range("B28").Value = 44192
$ws.Range("A27:B28").Font.Size = 11
$ws.Range("A27:B28").HorizontalAlignment = -4108
$ws.Range("A27:B28").VerticalAlignment = -4108

$ws.Range("C27:G27").Borders.Item(7).LineStyle = 1
$ws.Range("C27:G27").Borders.Item(9).LineStyle = 0
$ws.Range("G27").Borders.Item(10).LineStyle = 1
$ws.Range("C28:G28").Borders.Item(7).LineStyle = 1
$ws.Range("C28:G28").Borders.Item(9).LineStyle = 1
$ws.Range("G28").Borders.Item(10).LineStyle = 1

$ws.Range("C26:G28").Merge()
$ws.Rows.Item(27).RowHeight = 25.95
$ws.Rows.Item(28).RowHeight = 21

# --- 6. Rows 29-30 (normal data rows again, same look as rows 15 & 17) ---
$ws.Range("A15:G15").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Range("A17:G17").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("A29").Value = 7
$ws.Range("B29").Value = 44193
$ws.Range("C29").Value = "QMVAR 2.0"
$ws.Range("D29").Value = "Design issue fixing"
$ws.Range("F29").Value = "WIP"

$ws.Range("A30").Value = 8
$ws.Range("B30").Value = 44194
$ws.Range("C30").Value = "QMVAR 2.0"
$ws.Range("D30").Value = "Design issue fixing"
$ws.Range("F30").Value = "Completed"

# --- 7. Rows 31-33: blank spacer rows (touched, no content) ---
$ws.Range("A31:A33").Font.Size = 11
$ws.Range("D31:G33").Font.Size = 11

# --- 8. Update sheet view / selection ---
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("F33").Select()
